# "close to end of auction" -- append new draft picks to the draftpicks sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("draftpicks")

$picks = @{
    453 = @{Team='rippe'; Player='Eduardo Nunez'; Salary=1; Pos='CI'}
    454 = @{Team='deano'; Player='Scott Schebler'; Salary=1; Pos='OF'}
    455 = @{Team='bears'; Player='Mike Minor'; Salary=1; Pos='P'}
    456 = @{Team='marmaduke'; Player='Jace Fry'; Salary=1; Pos='P'}
    457 = @{Team='pasadena'; Player='Jay Bruce'; Salary=3; Pos='OF'}
    458 = @{Team='dembums'; Player='Matt Kemp'; Salary=2; Pos='DH'}
    459 = @{Team='allrise'; Player='Elias Diaz'; Salary=1; Pos='C'}
    460 = @{Team='dsb'; Player='Nick Hundley'; Salary=1; Pos='C'}
    461 = @{Team='pasadena'; Player='Mitch Keller'; Salary=3; Pos='P'}
    462 = @{Team='deano'; Player='Matt Duffy'; Salary=1; Pos='3B'}
    463 = @{Team='marmaduke'; Player='Raimel Tapia'; Salary=1; Pos='DH'}
    464 = @{Team='pasadena'; Player='Orlando Arcia'; Salary=1; Pos='MI'}
    465 = @{Team='allrise'; Player='Jeff McNeil'; Salary=2; Pos='DH'}
    466 = @{Team='allrise'; Player='Bradley Zimmer'; Salary=1; Pos='OF'}
    467 = @{Team='rippe'; Player='Mike Fiers'; Salary=1; Pos='P'}
    468 = @{Team='deano'; Player='Todd Frazier'; Salary=1; Pos='CI'}
    469 = @{Team='marmaduke'; Player='Fernando Romero'; Salary=1; Pos='P'}
    470 = @{Team='pasadena'; Player='Greg Allen'; Salary=1; Pos='OF'}
}

# 1) lay down all 18 rows: copy formatting down from the last existing row,
#    and fill in team/salary/position/date -- player name (col B) filled in separately below
$prevRow = 452
$ws.Range("A452:E452").Copy($ws.Range("A453:E453"))
$ws.Range("A453").Value = $picks[453].Team
$ws.Range("C453").Value = $picks[453].Salary
$ws.Range("D453").Value = $picks[453].Pos
$ws.Range("A453:E453").Copy($ws.Range("A454:E454"))
$ws.Range("A454").Value = $picks[454].Team
$ws.Range("C454").Value = $picks[454].Salary
$ws.Range("D454").Value = $picks[454].Pos
$ws.Range("A454:E454").Copy($ws.Range("A455:E455"))
$ws.Range("A455").Value = $picks[455].Team
$ws.Range("C455").Value = $picks[455].Salary
$ws.Range("D455").Value = $picks[455].Pos
$ws.Range("A455:E455").Copy($ws.Range("A456:E456"))
$ws.Range("A456").Value = $picks[456].Team
$ws.Range("C456").Value = $picks[456].Salary
$ws.Range("D456").Value = $picks[456].Pos
$ws.Range("A456:E456").Copy($ws.Range("A457:E457"))
$ws.Range("A457").Value = $picks[457].Team
$ws.Range("C457").Value = $picks[457].Salary
$ws.Range("D457").Value = $picks[457].Pos
$ws.Range("A457:E457").Copy($ws.Range("A458:E458"))
$ws.Range("A458").Value = $picks[458].Team
$ws.Range("C458").Value = $picks[458].Salary
$ws.Range("D458").Value = $picks[458].Pos
$ws.Range("A458:E458").Copy($ws.Range("A459:E459"))
$ws.Range("A459").Value = $picks[459].Team
$ws.Range("C459").Value = $picks[459].Salary
$ws.Range("D459").Value = $picks[459].Pos
$ws.Range("A459:E459").Copy($ws.Range("A460:E460"))
$ws.Range("A460").Value = $picks[460].Team
$ws.Range("C460").Value = $picks[460].Salary
$ws.Range("D460").Value = $picks[460].Pos
$ws.Range("A460:E460").Copy($ws.Range("A461:E461"))
$ws.Range("A461").Value = $picks[461].Team
$ws.Range("C461").Value = $picks[461].Salary
$ws.Range("D461").Value = $picks[461].Pos
$ws.Range("A461:E461").Copy($ws.Range("A462:E462"))
$ws.Range("A462").Value = $picks[462].Team
$ws.Range("C462").Value = $picks[462].Salary
$ws.Range("D462").Value = $picks[462].Pos
$ws.Range("A462:E462").Copy($ws.Range("A463:E463"))
$ws.Range("A463").Value = $picks[463].Team
$ws.Range("C463").Value = $picks[463].Salary
$ws.Range("D463").Value = $picks[463].Pos
$ws.Range("A463:E463").Copy($ws.Range("A464:E464"))
$ws.Range("A464").Value = $picks[464].Team
$ws.Range("C464").Value = $picks[464].Salary
$ws.Range("D464").Value = $picks[464].Pos
$ws.Range("A464:E464").Copy($ws.Range("A465:E465"))
$ws.Range("A465").Value = $picks[465].Team
$ws.Range("C465").Value = $picks[465].Salary
$ws.Range("D465").Value = $picks[465].Pos
$ws.Range("A465:E465").Copy($ws.Range("A466:E466"))
$ws.Range("A466").Value = $picks[466].Team
$ws.Range("C466").Value = $picks[466].Salary
$ws.Range("D466").Value = $picks[466].Pos
$ws.Range("A466:E466").Copy($ws.Range("A467:E467"))
$ws.Range("A467").Value = $picks[467].Team
$ws.Range("C467").Value = $picks[467].Salary
$ws.Range("D467").Value = $picks[467].Pos
$ws.Range("A467:E467").Copy($ws.Range("A468:E468"))
$ws.Range("A468").Value = $picks[468].Team
$ws.Range("C468").Value = $picks[468].Salary
$ws.Range("D468").Value = $picks[468].Pos
$ws.Range("A468:E468").Copy($ws.Range("A469:E469"))
$ws.Range("A469").Value = $picks[469].Team
$ws.Range("C469").Value = $picks[469].Salary
$ws.Range("D469").Value = $picks[469].Pos
$ws.Range("A469:E469").Copy($ws.Range("A470:E470"))
$ws.Range("A470").Value = $picks[470].Team
$ws.Range("C470").Value = $picks[470].Salary
$ws.Range("D470").Value = $picks[470].Pos

# 2) fill in player names (col B) in the order they were drafted/recorded
$ws.Range("B453").Value = $picks[453].Player
$ws.Range("B454").Value = $picks[454].Player
$ws.Range("B455").Value = $picks[455].Player
$ws.Range("B456").Value = $picks[456].Player
$ws.Range("B457").Value = $picks[457].Player
$ws.Range("B458").Value = $picks[458].Player
$ws.Range("B460").Value = $picks[460].Player
$ws.Range("B461").Value = $picks[461].Player
$ws.Range("B463").Value = $picks[463].Player
$ws.Range("B464").Value = $picks[464].Player
$ws.Range("B465").Value = $picks[465].Player
$ws.Range("B466").Value = $picks[466].Player
$ws.Range("B467").Value = $picks[467].Player
$ws.Range("B468").Value = $picks[468].Player
$ws.Range("B469").Value = $picks[469].Player
$ws.Range("B470").Value = $picks[470].Player
$ws.Range("B462").Value = $picks[462].Player
$ws.Range("B459").Value = $picks[459].Player

$ws.Range("D461").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 446
$win.ScrollColumn = 1
